$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")
$ws.Range("Z1").Value2 = "test(a,b,c)"
Write-Host $ws.Range("Z1").Value2
